$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The order of writes below is chosen to reproduce the shared-string table
# order observed in the target workbook (Excel assigns new shared-string
# indices in the order a distinct string is first written to a cell).

# Row 16: "Spacing" header
$ws.Range("A16").Value = "Spacing"

# Numeric inputs / labels, row 19-21 first (paper width filled in later)
$ws.Range("A19").Value = "Number columns"
$ws.Range("D19").Value = 4

$ws.Range("A20").Value = "Width of left and right margins"
$ws.Range("D20").Value = 6.35

$ws.Range("A21").Value = "Width of middle gutters"
$ws.Range("D21").Value = 6

$ws.Range("A23").Value = "Guide positions"

# Guide position value labels (RHS/LHS for columns 1-3, skipping "Col 1, LHS" for now)
$ws.Range("B25").Value = "Col 1, RHS"
$ws.Range("B26").Value = "Col 2, LHS"
$ws.Range("B27").Value = "Col 2, RHS"
$ws.Range("B28").Value = "Col 3, LHS"

# Back-fill paper width and the derived column-width row
$ws.Range("A18").Value = "Width of paper"
$ws.Range("D18").Value = 210

$ws.Range("A22").Value = "Therefore, width of columns"
$ws.Range("D22").Formula = "= (D18 - (2*D20) - (D19-1) * D21) / D19"

# Remaining guide position labels / formulas
$ws.Range("B29").Value = "Col 3, RHS"
$ws.Range("B30").Value = "Col 4, LHS"
$ws.Range("B31").Value = "Col 4, RHS"

$ws.Range("B24").Value = "Col 1, LHS"

# Formulas for the guide-position column
$ws.Range("C24").Formula = "=D20"
$ws.Range("C25").Formula = "=C24+D22"
$ws.Range("C26").Formula = "=C25+D21"
$ws.Range("C27").Formula = "=C26+D22"
$ws.Range("C28").Formula = "=C27+D21"
$ws.Range("C29").Formula = "=C28+D22"
$ws.Range("C30").Formula = "=C29+D21"
$ws.Range("C31").Formula = "=C30+D22"

# Update the view: scroll so A4 is top-left, and select C29
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("C29").Select()
